# Mirandola.xlsx update: "aggiornamento fino a 28/06 incluso"
# Appends daily COVID rows 270-301 (2021-05-28 .. 2021-06-28) to Sheet1,
# extending the dimension from A1:D269 to A1:D301.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the existing formatting (notably the date number-format style
# carried on column A) from the last populated row down across the new
# rows before writing values into them.
$ws.Range("A269:D269").Copy()
$ws.Range("A270:D301").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# columns: row, A (date serial), B (nuovi pos.), C (somma mobile 7gg.),
# D (somma mobile 7gg. per 100mila abitanti)
$newData = @(
  @(270, 44344, 2, 11, 45.60152557831026),
  @(271, 44345, 2, 12, 49.7471188127021),
  @(272, 44346, 3, 15, 62.18389851587763),
  @(273, 44347, 1, 12, 49.7471188127021),
  @(274, 44348, 2, 11, 45.60152557831026),
  @(275, 44349, 0, 11, 45.60152557831026),
  @(276, 44350, 10, 20, 82.91186468783683),
  @(277, 44351, 1, 19, 78.76627145344499),
  @(278, 44352, 0, 17, 70.47508498466131),
  @(279, 44353, 0, 14, 58.03830528148578),
  @(280, 44354, 1, 14, 58.03830528148578),
  @(281, 44355, 0, 12, 49.7471188127021),
  @(282, 44356, 0, 12, 49.7471188127021),
  @(283, 44357, 0, 2, 8.291186468783684),
  @(284, 44358, 1, 2, 8.291186468783684),
  @(285, 44359, 0, 2, 8.291186468783684),
  @(286, 44360, 0, 2, 8.291186468783684),
  @(287, 44361, 0, 1, 4.145593234391842),
  @(288, 44362, 0, 1, 4.145593234391842),
  @(289, 44363, 0, 1, 4.145593234391842),
  @(290, 44364, 0, 1, 4.145593234391842),
  @(291, 44365, 0, 0, 0),
  @(292, 44366, 0, 0, 0),
  @(293, 44367, 0, 0, 0),
  @(294, 44368, 0, 0, 0),
  @(295, 44369, 0, 0, 0),
  @(296, 44370, 0, 0, 0),
  @(297, 44371, 1, 1, 4.145593234391842),
  @(298, 44372, 0, 1, 4.145593234391842),
  @(299, 44373, 0, 1, 4.145593234391842),
  @(300, 44374, 0, 1, 4.145593234391842),
  @(301, 44375, 1, 2, 8.291186468783684)
)

foreach ($row in $newData) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = $row[1]
  $ws.Cells.Item($r, 2).Value = $row[2]
  $ws.Cells.Item($r, 3).Value = $row[3]
  $ws.Cells.Item($r, 4).Value = $row[4]
}
